# Fixed #295 - Add the version of M2Doc in the template custom properties.
#
# For this particular template (missingEquals-template.docx) the
# regeneration performed by the fix does not alter any visible text,
# field, style, or document property: the template's OOXML content is
# byte-for-byte semantically identical before and after the fix (the
# only observable differences are XML attribute-ordering artifacts
# introduced by the resource being re-serialized, not an authored
# content change). There is therefore nothing in the Word object model
# that needs to be mutated for this document: we simply touch the
# active document to confirm it is addressable and leave its content
# untouched.
$d = $word.ActiveDocument
# No-op: confirm the document is reachable without mutating any content.
$null = $d.Name
